$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the calibration / semicalibration values
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 20
$ws.Range("H2").Value = 0.9

$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 0.9

# Update the current selection on the sheet
$ws.Range("F3").Select()
